# Updates the cryptocurrency price/volume table to the latest scraped
# values (GitHub Actions refresh). A handful of rows also swap which
# coin occupies them (ranking reshuffle), which touches B/C/D/E together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.042.04'
$ws.Range('E2').Value = '  +4.46%  '
$ws.Range('D3').Value = '2.287.81'
$ws.Range('E3').Value = '  +4.92%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.635'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.05'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +9.68%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.667'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +13.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0982'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '59.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.105'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').Value = '2.626.47'
$ws.Range('E15').Value = '  +4.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.898'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.14%  '
$ws.Range('D18').Value = '2.279.05'
$ws.Range('E18').Value = '  +4.87%  '
$ws.Range('D19').Value = '42.941.74'
$ws.Range('E19').Value = '  +4.40%  '
$ws.Range('E20').Value = '  +6.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.16%  '
$ws.Range('E24').Value = '  +7.79%  '
$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.04%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.85%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.69'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  +5.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.36'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +11.86%  '
$ws.Range('E34').Value = '  +6.82%  '
$ws.Range('E35').Value = '  +8.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.49'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +28.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.85'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +22.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.126'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.81'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.16%  '
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.75%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +17.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.81%  '
$ws.Range('E44').Value = '  +12.69%  '
$ws.Range('E45').Value = '  +7.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '61.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('E49').Value = '  +5.04%  '
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '97.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.01%  '
